$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.441.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.638.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.05"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.174"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.637.52"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.167"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.05"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000190"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.79%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.121.04"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "72.289.87"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.66"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.635.50"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.99"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "379.91"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.92"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.36"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.55%  "

$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.41"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.84%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.772.32"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.94%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0959"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.16"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "523.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.44%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.36"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.33"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.114"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.35%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.66"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.08"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.79"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.71"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.545"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.99%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0263"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.44%  "
